$wb = $excel.ActiveWorkbook

# --- Separate "dividends" vs "tax withholding" labels on the Tax Withholding sheet ---
$wsTax = $wb.Worksheets.Item("Tax Withholding")
$wsTax.Range("B2").Value = "Tax Withholding (NVDA)"
$wsTax.Range("B3").Value = "Tax Withholding (APPL)"

# Column B no longer needs to fit the long "Withheld Tax on Dividends (...)" text,
# so narrow it back down to fit the shorter "Tax Withholding (...)" labels.
$wsTax.Columns.Item(2).ColumnWidth = 21.6

# --- Updated Foreign Currencies figures (fix a few AWV-related numbers) ---
$wsFx = $wb.Worksheets.Item("Foreign Currencies")
$wsFx.Range("B2").Value = 1217.91
$wsFx.Range("B3").Value = 100
$wsFx.Range("B4").Value = 100
